$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "DANH SACH NGUOI DUNG" to "Sheet1"
$ws.Name = "Sheet1"

# Drop the "Lien ket avatar" (avatar link) column (column E); the old
# "Bi xoa" column (F) shifts left into its place.
$ws.Columns.Item(5).EntireColumn.Delete()

# Resize the remaining columns to the new layout:
#  - B:C (Ten nguoi dung / Email) shrink from ~50.8 to ~36.8 chars
#  - D:E (Vai tro / Bi xoa) both end up at ~17.2 chars
$ws.Columns.Item(2).ColumnWidth = 36
$ws.Columns.Item(3).ColumnWidth = 36
$ws.Columns.Item(4).ColumnWidth = 16.571428571428573
$ws.Columns.Item(5).ColumnWidth = 16.571428571428573

# Move the active selection to D7, matching the saved view state
$ws.Range("D7").Select() | Out-Null
